$wb = $excel.ActiveWorkbook

# sheet1 "展览" -> Worksheets.Item(1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1376
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F14").Value = 141
$ws.Range("F15").Value = 90
$ws.Range("F17").Value = 178
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 671
$ws.Range("F23").Value = 45
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F29").Value = 155
$ws.Range("F30").Value = 5219
$ws.Range("F31").Value = 549
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 139
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 82
$ws.Range("F38").Value = 49
$ws.Range("F39").Value = 3
$ws.Range("F40").Value = 51
$ws.Range("F41").Value = 247
$ws.Range("F43").Value = 3965
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 94

# sheet2 "演出" -> Worksheets.Item(2)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 2

# sheet4 "全部类型" -> Worksheets.Item(4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 202
$ws.Range("F3").Value = 1376
$ws.Range("F4").Value = 19349
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7326
$ws.Range("F11").Value = 713
$ws.Range("F12").Value = 238
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 346
$ws.Range("F20").Value = 65
$ws.Range("F21").Value = 671
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 54
$ws.Range("F25").Value = 300
$ws.Range("F26").Value = 1058
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 5219
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 32
$ws.Range("F38").Value = 12419
$ws.Range("F39").Value = 1312
$ws.Range("F40").Value = 49
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 247
$ws.Range("F44").Value = 327
$ws.Range("F45").Value = 3965
